# Weekly update: insert two new price-report rows at the top of the
# date-sorted "Choclo" data block (rows 373-398 shift down to 375-400),
# carrying in the latest week's observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the existing row 373; this pushes the old
# rows 373:398 down to 375:400 and copies formatting from the row above
# (so the date column keeps its date/time number format).
$ws.Rows("373:374").Insert()

# New row 373
$ws.Cells.Item(373, 1).Value = 8
$ws.Cells.Item(373, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(373, 3).Value = 'Coquimbo'
$ws.Cells.Item(373, 4).Value = 44516
$ws.Cells.Item(373, 5).Value = 4
$ws.Cells.Item(373, 6).Value = 100112024
$ws.Cells.Item(373, 7).Value = 'Choclo'
$ws.Cells.Item(373, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(373, 9).Value = 'Primera'
$ws.Cells.Item(373, 10).Value = 400
$ws.Cells.Item(373, 11).Value = 28000
$ws.Cells.Item(373, 12).Value = 29000
$ws.Cells.Item(373, 13).Value = 28500
$ws.Cells.Item(373, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(373, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(373, 16).Value = 407
$ws.Cells.Item(373, 17).Value = 70
$ws.Cells.Item(373, 18).Value = 'Hortaliza'

# New row 374
$ws.Cells.Item(374, 1).Value = 8
$ws.Cells.Item(374, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(374, 3).Value = 'Coquimbo'
$ws.Cells.Item(374, 4).Value = 44516
$ws.Cells.Item(374, 5).Value = 4
$ws.Cells.Item(374, 6).Value = 100112024
$ws.Cells.Item(374, 7).Value = 'Choclo'
$ws.Cells.Item(374, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(374, 9).Value = 'Primera'
$ws.Cells.Item(374, 10).Value = 400
$ws.Cells.Item(374, 11).Value = 19000
$ws.Cells.Item(374, 12).Value = 20000
$ws.Cells.Item(374, 13).Value = 19500
$ws.Cells.Item(374, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(374, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(374, 16).Value = 279
$ws.Cells.Item(374, 17).Value = 70
$ws.Cells.Item(374, 18).Value = 'Hortaliza'

Write-Output "inserted rows 373-374; new dimension should be A1:R400"
